# Auto-generated cell updates for Sheets workbook (Jenova market-profit data)
# Applies per-cell value changes as captured in the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 994.2857
$ws.Range("I6").Value = 994.2857
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2982.8571
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -2870.8571
$ws.Range("N6").ClearContents()
$ws.Range("H9").Value = 14650.25
$ws.Range("I9").Value = 19333.666
$ws.Range("J9").Value = 600
$ws.Range("K9").Value = 19333.666
$ws.Range("L9").Value = 600
$ws.Range("M9").Value = -19164.666
$ws.Range("N9").Value = -938
$ws.Range("H12").Value = 4543.5884
$ws.Range("I12").Value = 3238.2727
$ws.Range("J12").Value = 6936.6665
$ws.Range("K12").Value = 3238.2727
$ws.Range("L12").Value = 6936.6665
$ws.Range("M12").Value = -3068.2727
$ws.Range("N12").Value = -7276.6665
$ws.Range("H18").Value = 920
$ws.Range("I18").Value = 400
$ws.Range("K18").Value = 400
$ws.Range("M18").Value = -116
$ws.Range("H29").Value = 1559.4
$ws.Range("I29").Value = 1450.5
$ws.Range("J29").Value = 1995
$ws.Range("K29").Value = 4351.5
$ws.Range("L29").Value = 5985
$ws.Range("M29").Value = -4070.5
$ws.Range("N29").Value = -6547
$ws.Range("H51").Value = 5347.625
$ws.Range("J51").Value = 7582
$ws.Range("L51").Value = 7582
$ws.Range("N51").Value = -8550
$ws.Range("H55").Value = 59674
$ws.Range("I55").Value = 210
$ws.Range("J55").Value = 77970.62
$ws.Range("K55").Value = 210
$ws.Range("L55").Value = 77970.62
$ws.Range("M55").Value = 4
$ws.Range("N55").Value = -78398.62
$ws.Range("H64").Value = 9201.2
$ws.Range("I64").Value = 5003
$ws.Range("K64").Value = 5003
$ws.Range("M64").Value = -4755
$ws.Range("H67").Value = 9201.2
$ws.Range("I67").Value = 5003
$ws.Range("K67").Value = 5003
$ws.Range("M67").Value = -4145
$ws.Range("H76").Value = 143009390
$ws.Range("I76").Value = 177615.17
$ws.Range("K76").Value = 177615.17
$ws.Range("M76").Value = -177300.17
$ws.Range("H79").Value = 143009390
$ws.Range("I79").Value = 177615.17
$ws.Range("K79").Value = 177615.17
$ws.Range("M79").Value = -176523.17
$ws.Range("H92").Value = 561.5
$ws.Range("I92").Value = 113.416664
$ws.Range("K92").Value = 113.416664
$ws.Range("M92").Value = 1134.583336
$ws.Range("H101").Value = 591.6667
$ws.Range("I101").Value = 583.3333
$ws.Range("J101").Value = 600
$ws.Range("K101").Value = 1749.9999
$ws.Range("L101").Value = 1800
$ws.Range("M101").Value = -127.9999
$ws.Range("N101").Value = -5044
$ws.Range("H107").Value = 203580.8
$ws.Range("I107").Value = 253251
$ws.Range("J107").Value = 4900
$ws.Range("K107").Value = 253251
$ws.Range("L107").Value = 4900
$ws.Range("M107").Value = -251331
$ws.Range("N107").Value = -8740
$ws.Range("H113").Value = 14531.125
$ws.Range("J113").Value = 17285.715
$ws.Range("L113").Value = 17285.715
$ws.Range("N113").Value = -23793.715
$ws.Range("H131").Value = 2963.8235
$ws.Range("I131").Value = 1909.5454
$ws.Range("J131").Value = 4896.6665
$ws.Range("K131").Value = 5728.6362
$ws.Range("L131").Value = 14689.9995
$ws.Range("M131").Value = -688.6361999999999
$ws.Range("N131").Value = -24769.9995
$ws.Range("H138").Value = 5343.116
$ws.Range("I138").Value = 1827.9166
$ws.Range("J138").Value = 6703.839
$ws.Range("K138").Value = 5483.7498
$ws.Range("L138").Value = 20111.517
$ws.Range("M138").Value = -343.7497999999996
$ws.Range("N138").Value = -30391.517

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1457
$ws.Range("J5").Value = 4849.5
$ws.Range("L5").Value = 4849.5
$ws.Range("N5").Value = -5073.5
$ws.Range("H61").Value = 3208.3142
$ws.Range("I61").Value = 1511.32
$ws.Range("K61").Value = 1511.32
$ws.Range("M61").Value = -1299.32
$ws.Range("H122").Value = 6559.6
$ws.Range("I122").Value = 5599.3335
$ws.Range("K122").Value = 16798.0005
$ws.Range("M122").Value = -14348.0005
$ws.Range("H132").Value = 4530.643
$ws.Range("I132").Value = 603.2222
$ws.Range("K132").Value = 1809.6666
$ws.Range("M132").Value = 720.3334
$ws.Range("H136").Value = 3208.3142
$ws.Range("I136").Value = 1511.32
$ws.Range("K136").Value = 4533.96
$ws.Range("M136").Value = -1983.96

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1457
$ws.Range("J4").Value = 4849.5
$ws.Range("L4").Value = 4849.5
$ws.Range("N4").Value = -5079.5
$ws.Range("H107").Value = 977.5
$ws.Range("I107").Value = 977.5
$ws.Range("K107").Value = 977.5
$ws.Range("M107").Value = 942.5
$ws.Range("H134").Value = 4684.25
$ws.Range("I134").Value = 1368.5
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 4105.5
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -1570.5
$ws.Range("N134").Value = -29070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2750.7446
$ws.Range("I31").Value = 2442.889
$ws.Range("J31").Value = 2941.8276
$ws.Range("K31").Value = 2442.889
$ws.Range("L31").Value = 2941.8276
$ws.Range("M31").Value = -2147.889
$ws.Range("N31").Value = -3531.8276
$ws.Range("H34").Value = 2750.7446
$ws.Range("I34").Value = 2442.889
$ws.Range("J34").Value = 2941.8276
$ws.Range("K34").Value = 2442.889
$ws.Range("L34").Value = 2941.8276
$ws.Range("M34").Value = -2240.889
$ws.Range("N34").Value = -3345.8276
$ws.Range("H94").Value = 1509.55
$ws.Range("I94").Value = 1213
$ws.Range("J94").Value = 1669.2307
$ws.Range("K94").Value = 1213
$ws.Range("L94").Value = 1669.2307
$ws.Range("M94").Value = -762
$ws.Range("N94").Value = -2571.2307
$ws.Range("H122").Value = 3023.4614
$ws.Range("I122").Value = 1793
$ws.Range("J122").Value = 4459
$ws.Range("K122").Value = 5379
$ws.Range("L122").Value = 13377
$ws.Range("M122").Value = -2929
$ws.Range("N122").Value = -18277
$ws.Range("H132").Value = 4447.2354
$ws.Range("J132").Value = 5499.4
$ws.Range("L132").Value = 16498.2
$ws.Range("N132").Value = -21558.2
$ws.Range("H134").Value = 4182.7334
$ws.Range("I134").Value = 3329.45
$ws.Range("K134").Value = 9988.349999999999
$ws.Range("M134").Value = -7453.349999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 271.9
$ws.Range("J23").Value = 316.4
$ws.Range("L23").Value = 949.1999999999999
$ws.Range("N23").Value = -1419.2
$ws.Range("H38").Value = 39.090908
$ws.Range("J38").Value = 36.5
$ws.Range("L38").Value = 109.5
$ws.Range("N38").Value = -803.5
$ws.Range("H117").Value = 1974.75
$ws.Range("J117").Value = 1999.6666
$ws.Range("L117").Value = 5998.9998
$ws.Range("N117").Value = -12882.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 418107.16
$ws.Range("I24").Value = 69000
$ws.Range("J24").Value = 460000
$ws.Range("K24").Value = 69000
$ws.Range("L24").Value = 460000
$ws.Range("M24").Value = -68827
$ws.Range("N24").Value = -460346
$ws.Range("H107").Value = 515
$ws.Range("I107").Value = 621.5
$ws.Range("J107").Value = 461.75
$ws.Range("K107").Value = 621.5
$ws.Range("L107").Value = 461.75
$ws.Range("M107").Value = 1298.5
$ws.Range("N107").Value = -4301.75
$ws.Range("H132").Value = 1254042.5
$ws.Range("I132").Value = 1503852
$ws.Range("K132").Value = 4511556
$ws.Range("M132").Value = -4509026

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1714.5333
$ws.Range("I22").Value = 725.3333
$ws.Range("K22").Value = 725.3333
$ws.Range("M22").Value = -430.3333
$ws.Range("H27").Value = 1714.5333
$ws.Range("I27").Value = 725.3333
$ws.Range("K27").Value = 725.3333
$ws.Range("M27").Value = -618.3333
$ws.Range("H46").Value = 4677.9565
$ws.Range("J46").Value = 6388
$ws.Range("L46").Value = 6388
$ws.Range("N46").Value = -6764
$ws.Range("H55").Value = 427.10257
$ws.Range("I55").Value = 306.5862
$ws.Range("J55").Value = 776.6
$ws.Range("K55").Value = 306.5862
$ws.Range("L55").Value = 776.6
$ws.Range("M55").Value = -133.5862
$ws.Range("N55").Value = -1122.6
$ws.Range("H68").Value = 136901.5
$ws.Range("I68").Value = 18044.6
$ws.Range("J68").Value = 334996.34
$ws.Range("K68").Value = 18044.6
$ws.Range("L68").Value = 334996.34
$ws.Range("M68").Value = -17295.6
$ws.Range("N68").Value = -336494.34
$ws.Range("H71").Value = 136901.5
$ws.Range("I71").Value = 18044.6
$ws.Range("J71").Value = 334996.34
$ws.Range("K71").Value = 90223
$ws.Range("L71").Value = 1674981.7
$ws.Range("M71").Value = -86479
$ws.Range("N71").Value = -1682469.7
$ws.Range("H122").Value = 1914819.1
$ws.Range("I122").Value = 5003502
$ws.Range("K122").Value = 15010506
$ws.Range("M122").Value = -15008056
$ws.Range("H132").Value = 5412.143
$ws.Range("I132").Value = 3377
$ws.Range("K132").Value = 10131
$ws.Range("M132").Value = -7601
$ws.Range("H136").Value = 4627.778
$ws.Range("I136").Value = 3072.4546
$ws.Range("K136").Value = 9217.3638
$ws.Range("M136").Value = -6667.363799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3169.7407
$ws.Range("I132").Value = 1804.15
$ws.Range("K132").Value = 5412.450000000001
$ws.Range("M132").Value = -2882.450000000001
$ws.Range("H135").Value = 49250
$ws.Range("J135").Value = 49250
$ws.Range("L135").Value = 49250
$ws.Range("N135").Value = -59390
$ws.Range("H136").Value = 288598.66
$ws.Range("I136").Value = 346888.22
$ws.Range("K136").Value = 1040664.66
$ws.Range("M136").Value = -1038114.66
